$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EPPlus - Some data")
$ws.Name = "EPPlus - Some data with header"

$ws.PageSetup.PrintArea = '$A$1:$F$12'
$ws.PageSetup.PrintTitleRows = '$1:$5'
